$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.385.74"
$ws.Range("E2").Value = "  -0.74%  "

$ws.Range("D3").Value = "1.572.77"
$ws.Range("E3").Value = "  -0.47%  "

$ws.Range("E4").Value = "  -0.15%  "

$c = $ws.Range("D5")
$c.Value = "'212.04"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.23%  "

$ws.Range("E6").Value = "  -0.45%  "

$ws.Range("E7").Value = "  -0.09%  "

$c = $ws.Range("D8")
$c.Value = "'44.32"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -5.50%  "

$c = $ws.Range("D9")
$c.Value = "'23.74"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.95%  "

$c = $ws.Range("D10")
$c.Value = "'0.246"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.80%  "

$c = $ws.Range("D11")
$c.Value = "'0.0587"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.84%  "

$ws.Range("E12").Value = "  +1.33%  "

$ws.Range("D13").Value = "1.795.01"
$ws.Range("E13").Value = "  -0.60%  "

$ws.Range("D14").Value = "1.571.69"
$ws.Range("E14").Value = "  +0.13%  "

$c = $ws.Range("D15")
$c.Value = "'3.69"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.52%  "

# Row 16/17: Polygon and WrappedBTC swap places (with updated price/volume)
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c = $ws.Range("D16")
$c.Value = "'0.516"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.36%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "28.355.68"
$ws.Range("E17").Value = "  -0.78%  "

$c = $ws.Range("D18")
$c.Value = "'61.54"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.29%  "

$c = $ws.Range("D19")
$c.Value = "'230.10"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.49%  "

$c = $ws.Range("D20")
$c.Value = "'7.42"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.21%  "

$ws.Range("D21").Value = "0.0₃0684"
$ws.Range("E21").Value = "  -1.51%  "

$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("E23").Value = "  +0.91%  "

$c = $ws.Range("D24")
$c.Value = "'9.04"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.33%  "

$c = $ws.Range("D25")
$c.Value = "'2.05"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.67%  "

$c = $ws.Range("D26")
$c.Value = "'151.14"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.07%  "

$c = $ws.Range("D27")
$c.Value = "'14.92"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.51%  "

$c = $ws.Range("D28")
$c.Value = "'6.36"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.63%  "

$ws.Range("E29").Value = "  -1.05%  "

$ws.Range("E30").Value = "  -0.13%  "

$c = $ws.Range("D31")
$c.Value = "'0.0482"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +3.83%  "

$ws.Range("E32").Value = "  -3.50%  "

$ws.Range("E33").Value = "  -0.82%  "

$c = $ws.Range("D34")
$c.Value = "'3.13"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.60%  "

$ws.Range("D35").Value = "1.385.15"
$ws.Range("E35").Value = "  -0.92%  "

$ws.Range("E36").Value = "  +5.28%  "

$c = $ws.Range("D37")
$c.Value = "'1.51"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.61%  "

$c = $ws.Range("D39")
$c.Value = "'2.64"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.34%  "

$ws.Range("E40").Value = "  -1.87%  "

$c = $ws.Range("D41")
$c.Value = "'0.519"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.52%  "

$ws.Range("E42").Value = "  -0.09%  "

$c = $ws.Range("D43")
$c.Value = "'1.90"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.16%  "

$c = $ws.Range("D44")
$c.Value = "'0.787"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.25%  "

$c = $ws.Range("D45")
$c.Value = "'0.0470"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.53%  "

$ws.Range("E46").Value = "  -4.39%  "

$c = $ws.Range("D47")
$c.Value = "'62.36"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.97%  "

$c = $ws.Range("D48")
$c.Value = "'0.921"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -5.86%  "

$ws.Range("D49").Value = "1.708.57"
$ws.Range("E49").Value = "  -0.48%  "

$ws.Range("E50").Value = "  +0.89%  "

$c = $ws.Range("D51")
$c.Value = "'85.15"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.04%  "
